$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 06:48"

# --- Row 6 (India): refreshed case counters ---
$ws.Range("B6").Value = 2836925
$ws.Range("C6").Value = 1103
$ws.Range("E6").Value = 686863

# --- Row 118: refreshed case counters ---
$ws.Range("B118").Value = 3389
$ws.Range("C118").Value = 7
$ws.Range("D118").Value = 3218
$ws.Range("E118").Value = 113

# --- Country reorder cluster: Jamaica moves above Togo / Niger ---
$ws.Range("A154").Value = "Jamaica"
$ws.Range("B154").Value = 1192
$ws.Range("C154").Value = 46
$ws.Range("D154").Value = 772
$ws.Range("E154").Value = 406
$ws.Range("H154").Value = 14

$ws.Range("A155").Value = "Togo"
$ws.Range("B155").Value = 1190
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 875
$ws.Range("E155").Value = 288
$ws.Range("H155").Value = 27

$ws.Range("A156").Value = "Niger"
$ws.Range("B156").Value = 1167
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 1079
$ws.Range("E156").Value = 19
$ws.Range("H156").Value = 69

# --- Country reorder cluster: Butan moves above Monaco ---
$ws.Range("A189").Value = "Butan"
$ws.Range("B189").Value = 150
$ws.Range("C189").Value = 3
$ws.Range("D189").Value = 105
$ws.Range("E189").Value = 45
$ws.Range("H189").Value = 0

$ws.Range("A190").Value = "Monaco"
$ws.Range("B190").Value = 148
$ws.Range("C190").Value = 0
$ws.Range("D190").Value = 114
$ws.Range("E190").Value = 30
$ws.Range("H190").Value = 4

# --- Country reorder cluster: Islas Malvinas moves above Montserrat ---
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
